$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title/timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 11:12"

# Update country names in column A that shifted position due to re-sorting
$ws.Range("A15").Value = "Belgica"
$ws.Range("A16").Value = "Austria"
$ws.Range("A94").Value = "Oman"
$ws.Range("A95").Value = "Venezuela"
$ws.Range("A97").Value = "Sri Lanka"
$ws.Range("A120").Value = "Banglades"
$ws.Range("A121").Value = "Ruanda"
$ws.Range("A122").Value = "Paraguay"
$ws.Range("A123").Value = "Puerto Rico"
$ws.Range("A133").Value = "Isla de Man"
$ws.Range("A134").Value = "Polinesia Francesa"
$ws.Range("A135").Value = "Guatemala"
$ws.Range("A143").Value = "El Salvador"
$ws.Range("A144").Value = "Tanzania"
$ws.Range("A146").Value = "Etiopia"
$ws.Range("A147").Value = "Zambia"
$ws.Range("A148").Value = "Dominica"
$ws.Range("A150").Value = "Republica de Yibuti"
$ws.Range("A153").Value = "Surinam"
$ws.Range("A154").Value = "Haiti"
$ws.Range("A156").Value = "Seychelles"
$ws.Range("A157").Value = "Bermudas"
$ws.Range("A158").Value = "Niger"
$ws.Range("A159").Value = "Gabon"
$ws.Range("A165").Value = "Siria"
$ws.Range("A166").Value = "Fiyi"
$ws.Range("A169").Value = "Suazilandia"
$ws.Range("A170").Value = "Santa Sede"
$ws.Range("A171").Value = "Guinea"
$ws.Range("A172").Value = "Eritrea"
$ws.Range("A175").Value = "Liberia"
$ws.Range("A176").Value = "San Martin (Parte Holandesa)"
$ws.Range("A178").Value = "San Bartolome"
$ws.Range("A180").Value = "Republica de Africa Central"
$ws.Range("A181").Value = "Republica del Chad"
$ws.Range("A182").Value = "Angola"
$ws.Range("A184").Value = "Santa Lucia"
$ws.Range("A185").Value = "Mauritania"
$ws.Range("A186").Value = "Gambia"
$ws.Range("A187").Value = "Sudan"
$ws.Range("A188").Value = "Nepal"
$ws.Range("A189").Value = "Zimbabue"
$ws.Range("A190").Value = "Belice"
$ws.Range("A191").Value = "Nicaragua"
$ws.Range("A192").Value = "Somalia"
$ws.Range("A193").Value = "Butan"
$ws.Range("A194").Value = "Islas Virgenes Britanicas"
$ws.Range("A195").Value = "San Cristobal y Nieves"
$ws.Range("A196").Value = "Guinea-Bisau"
$ws.Range("A197").Value = "Mali"
$ws.Range("A198").Value = "Islas Turcas y Caicos"
$ws.Range("A199").Value = "Montserrat"
$ws.Range("A201").Value = "Granada"
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "San Vicente y las Granadinas"
$ws.Range("A204").Value = "Libia"

# Update statistic values in columns B-H
$ws.Range("B8").Value = 39355
$ws.Range("C8").Value = 2032
$ws.Range("E8").Value = 35586
$ws.Range("B9").Value = 29406
$ws.Range("C9").Value = 2389
$ws.Range("D9").Value = 10457
$ws.Range("E9").Value = 16715
$ws.Range("G9").Value = 157
$ws.Range("H9").Value = 2234
$ws.Range("B15").Value = 6235
$ws.Range("C15").Value = 1298
$ws.Range("D15").Value = 547
$ws.Range("E15").Value = 5468
$ws.Range("F15").Value = 474
$ws.Range("G15").Value = 42
$ws.Range("H15").Value = 220
$ws.Range("B16").Value = 6001
$ws.Range("C16").Value = 413
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = 5950
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 42
$ws.Range("B18").Value = 3191
$ws.Range("C18").Value = 107
$ws.Range("E18").Value = 3171
$ws.Range("B33").Value = 1106
$ws.Range("C33").Value = 43
$ws.Range("E33").Value = 1077
$ws.Range("B93").Value = 114
$ws.Range("C93").Value = 5
$ws.Range("D93").Value = 5
$ws.Range("E93").Value = 109
$ws.Range("B94").Value = 109
$ws.Range("C94").Value = 10
$ws.Range("D94").Value = 23
$ws.Range("E94").Value = 86
$ws.Range("F94").Value = 0
$ws.Range("B95").Value = 106
$ws.Range("D95").Value = 15
$ws.Range("E95").Value = 91
$ws.Range("F95").Value = 2
$ws.Range("B96").Value = 105
$ws.Range("C96").Value = 6
$ws.Range("E96").Value = 96
$ws.Range("B97").Value = 102
$ws.Range("D97").Value = 3
$ws.Range("E97").Value = 99
$ws.Range("F97").Value = 3
$ws.Range("D108").Value = 1
$ws.Range("E108").Value = 63
$ws.Range("B120").Value = 44
$ws.Range("C120").Value = 5
$ws.Range("D120").Value = 11
$ws.Range("E120").Value = 28
$ws.Range("F120").Value = 1
$ws.Range("H120").Value = 5
$ws.Range("C121").Value = 0
$ws.Range("E121").Value = 41
$ws.Range("F121").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("B122").Value = 41
$ws.Range("C122").Value = 4
$ws.Range("D122").Value = 0
$ws.Range("E122").Value = 38
$ws.Range("F122").Value = 1
$ws.Range("H122").Value = 3
$ws.Range("D123").Value = 1
$ws.Range("E123").Value = 36
$ws.Range("H123").Value = 2
$ws.Range("C133").Value = 2
$ws.Range("B134").Value = 25
$ws.Range("D134").Value = 0
$ws.Range("E134").Value = 25
$ws.Range("H134").Value = 0
$ws.Range("B135").Value = 24
$ws.Range("D135").Value = 4
$ws.Range("E135").Value = 19
$ws.Range("H135").Value = 1
$ws.Range("C143").Value = 4
$ws.Range("C144").Value = 0
$ws.Range("C158").Value = 0
$ws.Range("C159").Value = 1
$ws.Range("C185").Value = 1
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 3
$ws.Range("D188").Value = 1
$ws.Range("H188").Value = 0
$ws.Range("B189").Value = 3
$ws.Range("H189").Value = 1
$ws.Range("C192").Value = 1
$ws.Range("C194").Value = 0
